$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (C1:J1) - multiply the old 1..8 sequence by 1.525
$ws.Range("C1").Value = 1.5249999999999999
$ws.Range("D1").Value = 3.05
$ws.Range("E1").Value = 4.5749999999999993
$ws.Range("F1").Value = 6.1
$ws.Range("G1").Value = 7.625
$ws.Range("H1").Value = 9.1499999999999986
$ws.Range("I1").Value = 10.674999999999999
$ws.Range("J1").Value = 12.2

# Update column A (A3:A16) - multiply the old 1..14 sequence by 1.525
$ws.Range("A3").Value = 1.5249999999999999
$ws.Range("A4").Value = 3.05
$ws.Range("A5").Value = 4.5749999999999993
$ws.Range("A6").Value = 6.1
$ws.Range("A7").Value = 7.625
$ws.Range("A8").Value = 9.1499999999999986
$ws.Range("A9").Value = 10.674999999999999
$ws.Range("A10").Value = 12.2
$ws.Range("A11").Value = 13.725
$ws.Range("A12").Value = 15.25
$ws.Range("A13").Value = 16.774999999999999
$ws.Range("A14").Value = 18.299999999999997
$ws.Range("A15").Value = 19.824999999999999
$ws.Range("A16").Value = 21.349999999999998

# Update the active cell selection to match the new view state
$ws.Range("I25").Select()
